$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "311.94"
Set-TextValue "E2" "1.01%"

# Row 3
Set-TextValue "D3" "38.03"
Set-TextValue "E3" "0.36%"

# Row 4
Set-TextValue "D4" "5.125"
Set-TextValue "E4" "0.91%"

# Row 5
Set-TextValue "D5" "0.07932"
Set-TextValue "E5" "0.31%"

# Row 6
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.402"
Set-TextValue "E6" "1.39%"

# Row 7
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.913"
Set-TextValue "E7" "-2.77%"

# Row 8
Set-TextValue "B8" "KuCoinToken"
Set-TextValue "C8" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "8.242"
Set-TextValue "E8" "-0.48%"

# Row 9
Set-TextValue "B9" "BTSEToken"
Set-TextValue "C9" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D9" "2.842"
Set-TextValue "E9" "-11.11%"

# Row 10
Set-TextValue "B10" "MXToken"
Set-TextValue "C10" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D10" "0.9265"
Set-TextValue "E10" "-0.71%"

# Row 11
Set-TextValue "B11" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.1204"
Set-TextValue "E11" "-7.80%"

# Row 12
Set-TextValue "B12" "WazirX"
Set-TextValue "C12" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D12" "0.1915"
Set-TextValue "E12" "-0.89%"

# Row 13
Set-TextValue "B13" "MandalaExchangeToken"
Set-TextValue "C13" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.09216"
Set-TextValue "E13" "4.68%"

# Row 14
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03396"
Set-TextValue "E14" "-0.62%"

# Row 15
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09625"
Set-TextValue "E15" "-1.16%"

# Row 16
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001372"
Set-TextValue "E16" "-1.01%"

# Row 17
Set-TextValue "B17" "TigerCash"
Set-TextValue "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.005850"
Set-TextValue "E17" "-0.88%"

# Row 18
Set-TextValue "B18" "LEO"
Set-TextValue "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.534"
Set-TextValue "E18" "-1.57%"

# Row 19
Set-TextValue "E19" "0.18%"

# Row 20
Set-TextValue "D20" "5.260"
Set-TextValue "E20" "5.21%"

# Row 21
Set-TextValue "E21" "-0.81%"

# Row 22
Set-TextValue "D22" "0.2587"
Set-TextValue "E22" "4.28%"

# Row 23
Set-TextValue "D23" "0.02100"
Set-TextValue "E23" "179.76%"

# Row 24
Set-TextValue "D24" "0.04367"
Set-TextValue "E24" "1.41%"

# Row 25
Set-TextValue "D25" "0.001248"
Set-TextValue "E25" "2.49%"

# Row 26
Set-TextValue "D26" "0.004280"
Set-TextValue "E26" "-7.38%"

# Row 27
Set-TextValue "D27" "0.0001297"
Set-TextValue "E27" "-63.89%"

# Row 39
Set-TextValue "D39" "0.02107"
Set-TextValue "E39" "-8.62%"

# Row 40
Set-TextValue "D40" "0.05095"
Set-TextValue "E40" "0.95%"

# Row 41
Set-TextValue "D41" "0.007636"
Set-TextValue "E41" "1.50%"

# Row 42
Set-TextValue "D42" "0.009128"
Set-TextValue "E42" "-7.91%"

# Row 43
Set-TextValue "D43" "0.1355"
Set-TextValue "E43" "-0.48%"

# Row 44
Set-TextValue "D44" "0.002033"
Set-TextValue "E44" "0.84%"

# Row 45
Set-TextValue "D45" "0.008645"
Set-TextValue "E45" "-1.44%"

# Row 46
Set-TextValue "D46" "0.00006676"
Set-TextValue "E46" "2.05%"

# Row 47
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "0.02%"

# Row 48
Set-TextValue "D48" "0.002894"
Set-TextValue "E48" "-3.38%"

# Row 49
Set-TextValue "D49" "0.001199"
Set-TextValue "E49" "-0.12%"

# Row 50
Set-TextValue "D50" "0.00002098"
Set-TextValue "E50" "0.02%"

# Row 51
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "0.02%"
